$wb = $excel.ActiveWorkbook

# --- Remove the two blank sheets (Sheet2, Sheet3) ---
$wsSheet2 = $wb.Worksheets.Item("Sheet2")
$wsSheet2.Delete()
$wsSheet3 = $wb.Worksheets.Item("Sheet3")
$wsSheet3.Delete()

# --- Rename the remaining sheet to the new board name ---
$ws = $wb.Worksheets.Item("Power Supply Board - Rev A")
$ws.Name = "Breakout Board - Rev B"

# --- Update the BOM title (row 1) ---
$ws.Range("A1").Value = "Bill of Materials for 'Marmote - Breakout Board Rev B (Smoky)'"

# --- Swap the two mezzanine-connector BOM rows (row 4 <-> row 5 details) ---
# Row 4 now carries the former row-5 ("TOP") part data
$ws.Range("C4").Value = "CON-71741-2184-TOP"
$ws.Range("D4").Value = "CON-71741-2184-TOP"
$ws.Range("F4").Value = "71741-0002"
$ws.Range("H4").Value = "WM3498CT-ND"
$ws.Range("I4").Value = "MEZZANINE 1MM BTB PLUG 84CKT"
$ws.Range("K4").Value = 6.24

# Row 5 now carries the former row-4 ("BOT") part data
$ws.Range("C5").Value = "CON-71741-2184-BOT"
$ws.Range("D5").Value = "CON-71741-2184-BOT"
$ws.Range("F5").Value = "71742-3003"
$ws.Range("H5").Value = "WM3499CT-ND"
$ws.Range("I5").Value = "MEZZANINE 1MM BTB RECPT 84CKT"
$ws.Range("K5").Value = 5.81

# --- Update the active selection on the sheet ---
$ws.Range("A1:L1").Select()

Write-Output "edit applied"
